# =====================================================================
# CCO_eCoaching_Dimension_Table_Data.xlsx - Internal version 29 update
# - Updates for TFS 549 and 644:
#   * CSR Survey: 4 new DIM tables (Survey_DIM_Type, Sirvey_DIM_Question,
#     Survey_DIM_Response, Survey_DIM_QAnswer)
#   * ARC Escalations OMR feeds: new Sub Coaching Reason record + CSR=0
#     for "OMR: Inappropriate ARC Escalation" in Coaching_Reason_Selection
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Revision_History: two new rows documenting this release
# ---------------------------------------------------------------------
$rev = $wb.Worksheets.Item("Revision_History")

$rev.Cells.Item(32,1).Value = 29
$rev.Cells.Item(32,2).Value = 42261
$rev.Cells.Item(2,2).Copy()
$rev.Cells.Item(32,2).PasteSpecial(-4122) | Out-Null
$rev.Cells.Item(32,3).Value = "Susmitha Palacherla"
$rev.Cells.Item(32,4).Value = 549
$rev.Cells.Item(32,5).Value = "CSR Survey. Added 4 new Tabs for CSR Survey DIM tables."

$rev.Cells.Item(33,1).Value = 29.1
$rev.Cells.Item(33,2).Value = 42268
$rev.Cells.Item(2,2).Copy()
$rev.Cells.Item(33,2).PasteSpecial(-4122) | Out-Null
$rev.Cells.Item(33,3).Value = "Susmitha Palacherla"
$rev.Cells.Item(33,4).Value = 644
$rev.Cells.Item(33,5).Value = "ARC Escalations OMR feeds.Updated SubCoachingReason tab to add new record and set CSR to 0 in Coaching_Reason_selection table for IAE."

# Row 21 height tweak (was 45, now 30)
$rev.Rows.Item(21).RowHeight = 30

# New selection on Revision_History (no longer the active tab)
$rev.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Coaching_Reason_Selection: set CSR=0 for "OMR: Inappropriate ARC
#    Escalation" sub coaching reason row (IAE)
# ---------------------------------------------------------------------
$crs = $wb.Worksheets.Item("Coaching_Reason_Selection")
$crs.Cells.Item(189,10).Value = 0

# ---------------------------------------------------------------------
# 3) DIM_Sub_Coaching_Reason: add new record
#    "OMR: Inappropriate ARC Transfer"
# ---------------------------------------------------------------------
$subReason = $wb.Worksheets.Item("DIM_Sub_Coaching_Reason")
$subReason.Cells.Item(233,1).Value = 231
$subReason.Cells.Item(233,2).Value = "OMR: Inappropriate ARC Transfer"

# ---------------------------------------------------------------------
# 4) CSR Survey: 4 new DIM tables
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# --- Survey_DIM_Type -----------------------------------------------
$surveyType = $wb.Worksheets.Add($null, $lastSheet)
$surveyType.Name = "Survey_DIM_Type"

$surveyType.Cells.Item(1,1).Value  = "SurveyTypeID"
$surveyType.Cells.Item(1,2).Value  = "Decsription"
$surveyType.Cells.Item(1,3).Value  = "StartDate"
$surveyType.Cells.Item(1,4).Value  = "EndDate"
$surveyType.Cells.Item(1,5).Value  = "isActive"
$surveyType.Cells.Item(1,6).Value  = "CSR"
$surveyType.Cells.Item(1,7).Value  = "Supervisor"
$surveyType.Cells.Item(1,8).Value  = "Quality"
$surveyType.Cells.Item(1,9).Value  = "LSA"
$surveyType.Cells.Item(1,10).Value = "Training"
$surveyType.Cells.Item(1,11).Value = "LastUpdateDate"

$surveyType.Cells.Item(2,1).Value  = 1
$surveyType.Cells.Item(2,2).Value  = "Employee Survey"
$surveyType.Cells.Item(2,3).Value  = 20150901
$surveyType.Cells.Item(2,4).Value  = 99991231
$surveyType.Cells.Item(2,5).Value  = 1
$surveyType.Cells.Item(2,6).Value  = 1
$surveyType.Cells.Item(2,7).Value  = 0
$surveyType.Cells.Item(2,8).Value  = 0
$surveyType.Cells.Item(2,9).Value  = 0
$surveyType.Cells.Item(2,10).Value = 0
$surveyType.Cells.Item(2,11).Value = 42248

$surveyType.Columns.Item(11).ColumnWidth = 23.29

# Mint the LastUpdateDate number format (yyyy-mm-dd hh:mm:ss:mss) on this
# range; subsequent sheets copy this exact style so they all share one
# cellXf.
$surveyType.Range("K1:K2").NumberFormat = "yyyy-mm-dd hh:mm:ss:mss"
$dateStyleSrc = $surveyType.Range("K2")

$surveyType.Range("A1:K1").Select() | Out-Null

# --- Sirvey_DIM_Question ---------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$surveyQ = $wb.Worksheets.Add($null, $lastSheet)
$surveyQ.Name = "Sirvey_DIM_Question"

$surveyQ.Cells.Item(1,1).Value = "QuestionID"
$surveyQ.Cells.Item(1,2).Value = "Description"
$surveyQ.Cells.Item(1,3).Value = "DisplayOrder"
$surveyQ.Cells.Item(1,4).Value = "StartDate"
$surveyQ.Cells.Item(1,5).Value = "EndDate"
$surveyQ.Cells.Item(1,6).Value = "isHotTopic"
$surveyQ.Cells.Item(1,7).Value = "isActive"
$surveyQ.Cells.Item(1,8).Value = "LastUpdateDate"

$questions = @(
    @(-1, "Unknown", 0, 0),
    @(1,  "Was the call played back for you during your last coaching session? (If applicable). |     If no, what reason was provided?", 1, 1),
    @(2,  "Will you be able to apply the information from your last coaching session? |    If yes, how?  If no, why  not?", 2, 1),
    @(3,  "Did you find the coaching session valuable/effective? |   If yes, what specifically.  If no, why not?", 3, 1),
    @(4,  "Please rate the effectiveness of the coaching notes provided in the eCL. |   Please explain below.", 4, 1),
    @(5,  "Please rate your overall coaching experience. |   Please explain below.", 5, 1)
)

$r = 2
foreach ($q in $questions) {
    $surveyQ.Cells.Item($r,1).Value = $q[0]
    $surveyQ.Cells.Item($r,2).Value = $q[1]
    $surveyQ.Cells.Item($r,3).Value = $q[2]
    $surveyQ.Cells.Item($r,4).Value = 20150901
    $surveyQ.Cells.Item($r,5).Value = 99991231
    $surveyQ.Cells.Item($r,6).Value = 0
    $surveyQ.Cells.Item($r,7).Value = $q[3]
    $surveyQ.Cells.Item($r,8).Value = 42248
    $r = $r + 1
}

$surveyQ.Columns.Item(1).ColumnWidth = 11
$surveyQ.Columns.Item(2).ColumnWidth = 119.29
$surveyQ.Columns.Item(3).ColumnWidth = 9.29
$surveyQ.Columns.Item(4).ColumnWidth = 9
$surveyQ.Columns.Item(5).ColumnWidth = 9
$surveyQ.Columns.Item(6).ColumnWidth = 9.29
$surveyQ.Columns.Item(7).ColumnWidth = 9.29
$surveyQ.Columns.Item(8).ColumnWidth = 22

$dateStyleSrc.Copy()
$surveyQ.Range("H1:H7").PasteSpecial(-4122) | Out-Null

$surveyQ.PageSetup.Orientation = 1

# --- Survey_DIM_Response ----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$surveyR = $wb.Worksheets.Add($null, $lastSheet)
$surveyR.Name = "Survey_DIM_Response"

$surveyR.Cells.Item(1,1).Value = "ResponseID"
$surveyR.Cells.Item(1,2).Value = "Value"
$surveyR.Cells.Item(1,3).Value = "isActive"
$surveyR.Cells.Item(1,4).Value = "LastUpdateDate"

$responses = @(
    @(-1, "Not Applicable"),
    @(1,  "Yes"),
    @(2,  "No"),
    @(3,  "N/A"),
    @(4,  "1 - Very Ineffective"),
    @(5,  "2 - Ineffective"),
    @(6,  "3 - Neither"),
    @(7,  "4 - Effective"),
    @(8,  "5 - Very Effective"),
    @(9,  "1 - Very Dissatisfied"),
    @(10, "2 - Dissatisfied"),
    @(11, "4 - Satisfied"),
    @(12, "5 - Very Satisfied")
)

$r = 2
foreach ($resp in $responses) {
    $surveyR.Cells.Item($r,1).Value = $resp[0]
    $surveyR.Cells.Item($r,2).Value = $resp[1]
    $surveyR.Cells.Item($r,3).Value = 1
    $surveyR.Cells.Item($r,4).Value = 42248
    $r = $r + 1
}

$surveyR.Columns.Item(1).ColumnWidth = 11.43
$surveyR.Columns.Item(2).ColumnWidth = 18.71
$surveyR.Columns.Item(3).ColumnWidth = 8

$dateStyleSrc.Copy()
$surveyR.Range("D1:D14").PasteSpecial(-4122) | Out-Null

# --- Survey_DIM_QAnswer ------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$surveyQA = $wb.Worksheets.Add($null, $lastSheet)
$surveyQA.Name = "Survey_DIM_QAnswer"

$surveyQA.Cells.Item(1,1).Value  = "SurveyTypeID"
$surveyQA.Cells.Item(1,2).Value  = "QuestionID"
$surveyQA.Cells.Item(1,3).Value  = "QuestionNumber"
$surveyQA.Cells.Item(1,4).Value  = "ResponseID"
$surveyQA.Cells.Item(1,5).Value  = "ResponseValue"
$surveyQA.Cells.Item(1,6).Value  = "ResponseType"
$surveyQA.Cells.Item(1,7).Value  = "isHotTopic"
$surveyQA.Cells.Item(1,8).Value  = "StartDate"
$surveyQA.Cells.Item(1,9).Value  = "EndDate"
$surveyQA.Cells.Item(1,10).Value = "isActive"
$surveyQA.Cells.Item(1,11).Value = "LastUpdateDate"

$qanswers = @(
    @(1,1,1,1,  "Yes"),
    @(1,1,1,2,  "No"),
    @(1,1,1,3,  "N/A"),
    @(1,2,2,1,  "Yes"),
    @(1,2,2,2,  "No"),
    @(1,3,3,1,  "Yes"),
    @(1,3,3,2,  "No"),
    @(1,4,4,4,  "1 - Very Ineffective"),
    @(1,4,4,5,  "2 - Ineffective"),
    @(1,4,4,6,  "3 - Neither"),
    @(1,4,4,7,  "4 - Effective"),
    @(1,4,4,8,  "5 - Very Effective"),
    @(1,5,5,9,  "1 - Very Dissatisfied"),
    @(1,5,5,10, "2 - Dissatisfied"),
    @(1,5,5,6,  "3 - Neither"),
    @(1,5,5,11, "4 - Satisfied"),
    @(1,5,5,12, "5 - Very Satisfied")
)

$r = 2
foreach ($qa in $qanswers) {
    $surveyQA.Cells.Item($r,1).Value  = $qa[0]
    $surveyQA.Cells.Item($r,2).Value  = $qa[1]
    $surveyQA.Cells.Item($r,3).Value  = $qa[2]
    $surveyQA.Cells.Item($r,4).Value  = $qa[3]
    $surveyQA.Cells.Item($r,5).Value  = $qa[4]
    $surveyQA.Cells.Item($r,6).Value  = "Radio Button"
    $surveyQA.Cells.Item($r,7).Value  = 0
    $surveyQA.Cells.Item($r,8).Value  = 20150901
    $surveyQA.Cells.Item($r,9).Value  = 99991231
    $surveyQA.Cells.Item($r,10).Value = 1
    $surveyQA.Cells.Item($r,11).Value = 42248
    $r = $r + 1
}

$surveyQA.Columns.Item(11).ColumnWidth = 24

$dateStyleSrc.Copy()
$surveyQA.Range("K1:K18").PasteSpecial(-4122) | Out-Null

$surveyQA.PageSetup.Orientation = 1

# Survey_DIM_QAnswer is the active sheet/tab when the workbook was saved
$surveyQA.Activate() | Out-Null
$surveyQA.Range("A1").Select() | Out-Null

Write-Output "done"
